$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K (strikeouts) column values for rows 2-31 (column G)
$newK = @{
    2  = 7
    3  = 10
    4  = 11
    5  = 10
    6  = 11
    7  = 5
    8  = 4
    9  = 10
    10 = 5
    11 = 8
    12 = 5
    13 = 3
    14 = 9
    15 = 8
    16 = 4
    17 = 9
    18 = 8
    19 = 10
    20 = 4
    21 = 6
    22 = 6
    23 = 6
    24 = 12
    25 = 6
    26 = 9
    27 = 2
    28 = 2
    29 = 5
    30 = 5
    31 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
